$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the formatting of the last existing row (14) down onto the new
# row (15) before writing values, so the new cells pick up the same style
# indices (s="1" for appid, s="2" for the two email-like columns, default
# for the rest) that the author's row already used.
$ws.Range("A14:G14").Copy()
$ws.Range("A15:G15").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New review row appended to the table.
$ws.Range("A15").Value = "com.singleton.strechy"
$ws.Range("B15").Value = "taxi"
$ws.Range("C15").Value = "stavsade45@gmail.com"
$ws.Range("D15").Value = "galiatia942@gmail.com"
$ws.Range("E15").Value = "27/5/2019 15:59"
$ws.Range("F15").Value = "this is the game! My game! And happy forever and after"
$ws.Range("G15").Value = "yes"

# Match the author's final selection/scroll position after the edit.
$ws.Range("G16").Select() | Out-Null
